$wb = $excel.ActiveWorkbook

$wsModify = $wb.Worksheets.Item("Modify Transaction")

# Insert a new row above the existing row 4 ("submitmakerepayment"/"click")
$wsModify.Rows.Item(4).Insert()

# Populate the newly inserted row with the ReceiptNumber label/value
$wsModify.Cells.Item(4, 1).Value = "ReceiptNumber"
$wsModify.Cells.Item(4, 2).Value = 8765432

# Update the selection on the "Modify Transaction" sheet
$wsModify.Range("D10").Select()

# Make "Modify Transaction" the active sheet (was "Transactions" before)
$wsModify.Activate()

$wb.Save()
